$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($null -ne $val -and $val -is [string] -and $val.EndsWith("16")) {
        $cell.Value = $val.Substring(0, $val.Length - 2)
    }
}
